$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H26").Value = 6505
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 6505
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 6505
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -7193
$ws.Range("H38").Value = 8751.762000000001
$ws.Range("J38").Value = 6332.778
$ws.Range("L38").Value = 18998.334
$ws.Range("N38").Value = -19742.334
$ws.Range("H74").Value = 4776.4
$ws.Range("I74").Value = 3480.5715
$ws.Range("K74").Value = 3480.5715
$ws.Range("M74").Value = -2544.5715
$ws.Range("H77").Value = 4776.4
$ws.Range("I77").Value = 3480.5715
$ws.Range("K77").Value = 17402.8575
$ws.Range("M77").Value = -12722.8575
$ws.Range("H92").Value = 775.41174
$ws.Range("I92").Value = 525.53845
$ws.Range("K92").Value = 525.53845
$ws.Range("M92").Value = 722.46155
$ws.Range("H106").Value = 4346.1665
$ws.Range("I106").Value = 3519.5
$ws.Range("K106").Value = 3519.5
$ws.Range("M106").Value = -2888.5
$ws.Range("H136").Value = 338963
$ws.Range("J136").Value = 338963
$ws.Range("L136").Value = 338963
$ws.Range("N136").Value = -349163
$ws.Range("H137").Value = 3624322.2
$ws.Range("J137").Value = 7247193.5
$ws.Range("L137").Value = 21741580.5
$ws.Range("N137").Value = -21746680.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2489
$ws.Range("I2").Value = 1054.1428
$ws.Range("K2").Value = 1054.1428
$ws.Range("M2").Value = -941.1428000000001
$ws.Range("H45").Value = 17053.428
$ws.Range("I45").Value = 21274.8
$ws.Range("K45").Value = 21274.8
$ws.Range("M45").Value = -20897.8
$ws.Range("H88").Value = 1956.6666
$ws.Range("J88").Value = 1952.4117
$ws.Range("L88").Value = 1952.4117
$ws.Range("N88").Value = -2764.4117
$ws.Range("H91").Value = 1956.6666
$ws.Range("J91").Value = 1952.4117
$ws.Range("L91").Value = 1952.4117
$ws.Range("N91").Value = -4760.411700000001
$ws.Range("H97").Value = 460
$ws.Range("I97").Value = 460
$ws.Range("K97").Value = 460
$ws.Range("M97").Value = 36
$ws.Range("H110").Value = 7454.5454
$ws.Range("I110").Value = 7454.5454
$ws.Range("K110").Value = 7454.5454
$ws.Range("M110").Value = -5409.5454
$ws.Range("H116").Value = 2489
$ws.Range("I116").Value = 1054.1428
$ws.Range("K116").Value = 1054.1428
$ws.Range("M116").Value = 1239.8572

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2489
$ws.Range("I3").Value = 1054.1428
$ws.Range("K3").Value = 1054.1428
$ws.Range("M3").Value = -940.1428000000001
$ws.Range("H86").Value = 45465390
$ws.Range("I86").Value = 18098.5
$ws.Range("J86").Value = 100002140
$ws.Range("K86").Value = 18098.5
$ws.Range("L86").Value = 100002140
$ws.Range("M86").Value = -16975.5
$ws.Range("N86").Value = -100004386
$ws.Range("H89").Value = 45465390
$ws.Range("I89").Value = 18098.5
$ws.Range("J89").Value = 100002140
$ws.Range("K89").Value = 90492.5
$ws.Range("L89").Value = 500010700
$ws.Range("M89").Value = -84876.5
$ws.Range("N89").Value = -500021932
$ws.Range("H99").Value = 4006.923
$ws.Range("I99").Value = 2787.7778
$ws.Range("J99").Value = 6750
$ws.Range("K99").Value = 2787.7778
$ws.Range("L99").Value = 6750
$ws.Range("M99").Value = -1289.7778
$ws.Range("N99").Value = -9746
$ws.Range("H107").Value = 3692.4783
$ws.Range("I107").Value = 840.53845
$ws.Range("K107").Value = 840.53845
$ws.Range("M107").Value = 1079.46155
$ws.Range("H134").Value = 626620.75
$ws.Range("I134").Value = 1601.697
$ws.Range("K134").Value = 4805.090999999999
$ws.Range("M134").Value = -2270.090999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 456.25
$ws.Range("J4").Value = 500
$ws.Range("L4").Value = 500
$ws.Range("N4").Value = -724
$ws.Range("H31").Value = 5404.171
$ws.Range("I31").Value = 1326.6666
$ws.Range("J31").Value = 8595.261
$ws.Range("K31").Value = 1326.6666
$ws.Range("L31").Value = 8595.261
$ws.Range("M31").Value = -1031.6666
$ws.Range("N31").Value = -9185.261
$ws.Range("H34").Value = 5404.171
$ws.Range("I34").Value = 1326.6666
$ws.Range("J34").Value = 8595.261
$ws.Range("K34").Value = 1326.6666
$ws.Range("L34").Value = 8595.261
$ws.Range("M34").Value = -1124.6666
$ws.Range("N34").Value = -8999.261
$ws.Range("H86").Value = 8474.684999999999
$ws.Range("I86").Value = 8749.5
$ws.Range("K86").Value = 8749.5
$ws.Range("M86").Value = -7626.5
$ws.Range("H89").Value = 8474.684999999999
$ws.Range("I89").Value = 8749.5
$ws.Range("K89").Value = 43747.5
$ws.Range("M89").Value = -38131.5
$ws.Range("H99").Value = 7524300
$ws.Range("I99").Value = 3705
$ws.Range("K99").Value = 3705
$ws.Range("M99").Value = -2207
$ws.Range("H107").Value = 506.7
$ws.Range("I107").Value = 516.3333
$ws.Range("K107").Value = 516.3333
$ws.Range("M107").Value = 1403.6667
$ws.Range("H126").Value = 7524300
$ws.Range("I126").Value = 3705
$ws.Range("K126").Value = 11115
$ws.Range("M126").Value = -8645
$ws.Range("H134").Value = 2165.1428
$ws.Range("I134").Value = 1848.4
$ws.Range("K134").Value = 5545.200000000001
$ws.Range("M134").Value = -3010.200000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 624.5
$ws.Range("J103").Value = 666.3333
$ws.Range("L103").Value = 1998.9999
$ws.Range("N103").Value = -3756.9999
$ws.Range("H107").Value = 1789.963
$ws.Range("J107").Value = 1903
$ws.Range("L107").Value = 5709
$ws.Range("N107").Value = -9549
$ws.Range("H113").Value = 1447.6666
$ws.Range("J113").Value = 996.5
$ws.Range("L113").Value = 2989.5
$ws.Range("N113").Value = -7329.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 166672460
$ws.Range("I80").Value = 400002720
$ws.Range("J80").Value = 7999.143
$ws.Range("K80").Value = 400002720
$ws.Range("L80").Value = 7999.143
$ws.Range("M80").Value = -400001722
$ws.Range("N80").Value = -9995.143
$ws.Range("H83").Value = 166672460
$ws.Range("I83").Value = 400002720
$ws.Range("J83").Value = 7999.143
$ws.Range("K83").Value = 2000013600
$ws.Range("L83").Value = 39995.715
$ws.Range("M83").Value = -2000008608
$ws.Range("N83").Value = -49979.715

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 7827.364
$ws.Range("I68").Value = 1501
$ws.Range("K68").Value = 1501
$ws.Range("M68").Value = -752
$ws.Range("H71").Value = 7827.364
$ws.Range("I71").Value = 1501
$ws.Range("K71").Value = 7505
$ws.Range("M71").Value = -3761
$ws.Range("H122").Value = 4484.533
$ws.Range("I122").Value = 4289.923
$ws.Range("J122").Value = 5749.5
$ws.Range("K122").Value = 12869.769
$ws.Range("L122").Value = 17248.5
$ws.Range("M122").Value = -10419.769
$ws.Range("N122").Value = -22148.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10651544
$ws.Range("I62").Value = 250001340
$ws.Range("J62").Value = 13775.267
$ws.Range("K62").Value = 250001340
$ws.Range("L62").Value = 13775.267
$ws.Range("M62").Value = -250000716
$ws.Range("N62").Value = -15023.267
$ws.Range("H65").Value = 10651544
$ws.Range("I65").Value = 250001340
$ws.Range("J65").Value = 13775.267
$ws.Range("K65").Value = 1250006700
$ws.Range("L65").Value = 68876.33499999999
$ws.Range("M65").Value = -1250003580
$ws.Range("N65").Value = -75116.33499999999
$ws.Range("H81").Value = 1852.8334
$ws.Range("I81").Value = 1823.4
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 3646.8
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -2585.8
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 1852.8334
$ws.Range("I84").Value = 1823.4
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 18234
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -12930
$ws.Range("N84").Value = -30608
$ws.Range("H96").Value = 3488.1428
$ws.Range("I96").Value = 1604.25
$ws.Range("J96").Value = 6000
$ws.Range("K96").Value = 1604.25
$ws.Range("L96").Value = 6000
$ws.Range("M96").Value = -231.25
$ws.Range("N96").Value = -8746
